# Weekly driver report update for 2025-04-19 (UCSantaCruz_driver_summary)
#
# The "Bad Drivers" table loses its worst entry (the Nineplus USB LAN
# adapter dropped out of the bad-roaming list), so its one remaining row
# (previously the MediaTek 24.34.2.571 entry on row 4) becomes row 3 with
# refreshed sample counts, and the Totals row below it is recomputed to
# match. Deleting the old row 4 shifts everything below it up by one row,
# which also carries the "Good Drivers" section up from rows 11-13 to
# rows 10-12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the Nineplus row entirely; rows 5+ shift up to fill the gap.
$ws.Rows.Item(4).Delete()

# Row 3 now holds what used to be the MediaTek 24.34.2.571 row; refresh it
# with this week's numbers.
$ws.Range("A3").Value = "MediaTek Wi-Fi 6 MT7921 Wireless LAN Card - 24.34.2.571"
$ws.Range("C3").Value = 941
$ws.Range("D3").Value = 55.3

# Totals row (now row 4) reflects the single remaining bad-driver entry.
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 941
